$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 / column D: append the new recipient to the existing Email list,
# then turn the cell into a mailto hyperlink (this is what picks up the
# built-in "Hyperlink" style automatically).
$newEmail = "Justine.Eisenhour@ssc-spc.gc.ca,ssc.ppcoesupport-cdeppsupport.spc@ssc-spc.gc.ca,najet.nouisser@ssc-spc.gc.ca"
$cell = $ws.Range("D4")
$cell.Value = $newEmail
$ws.Hyperlinks.Add($cell, "mailto:" + $newEmail) | Out-Null

# Widen column D to fit the long hyperlink text, and move the selection
# there, matching the saved view state.
$ws.Columns("D").ColumnWidth = 155.6667
$cell.Select() | Out-Null
